$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (values like "1.004" would otherwise
# be auto-parsed as numbers by Excel's smart type detection).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.996.63'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.638.78'
$ws.Range("E3").Value = '  -1.72%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '206.77'
$ws.Range("E5").Value = '  -1.61%  '
$ws.Range("D6").Value = '0.5125'
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '0.2556'
$ws.Range("E8").Value = '  -3.05%  '
$ws.Range("D9").Value = '0.06177'
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").Value = '20.46'
$ws.Range("E10").Value = '  -3.30%  '
$ws.Range("D11").Value = '0.07570'
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("D12").Value = '1.643.15'
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").Value = '4.351'
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("D14").Value = '1.862.02'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").Value = '0.5321'
$ws.Range("E15").Value = '  -4.96%  '
$ws.Range("D16").Value = '0.0₅7963'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '65.21'
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '26.011.70'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '4.624'
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("D21").Value = '185.27'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = '9.965'
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("D23").Value = '6.064'
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").Value = '146.94'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '7.329'
$ws.Range("E26").Value = '  -3.09%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.1193'
$ws.Range("E27").Value = '  -3.99%  '
$ws.Range("D28").Value = '15.46'
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("D29").Value = '1.350'
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("D30").Value = '0.06001'
$ws.Range("E30").Value = '  -4.25%  '
$ws.Range("D31").Value = '1.243'
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("D32").Value = '3.389'
$ws.Range("E32").Value = '  -2.41%  '
$ws.Range("D33").Value = '3.366'
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").Value = '1.610'
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("D35").Value = '0.9633'
$ws.Range("E35").Value = '  -3.31%  '
$ws.Range("D36").Value = '2.382'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").Value = '2.721'
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("D38").Value = '0.5811'
$ws.Range("E38").Value = '  -3.63%  '
$ws.Range("D39").Value = '0.01579'
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").Value = '1.074.01'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '5.801'
$ws.Range("E41").Value = '  -5.36%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8420'
$ws.Range("E43").Value = '  -2.46%  '
$ws.Range("D44").Value = '99.93'
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").Value = '1.795.12'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("D46").Value = '0.0₈105'
$ws.Range("E46").Value = '  -3.21%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").Value = '54.10'
$ws.Range("E48").Value = '  -3.42%  '
$ws.Range("D49").Value = '7.972'
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").Value = '0.05208'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  -0.28%  '

# Restore the default (no explicit) style on the Price column so the
# saved file does not carry a lingering text-number-format style.
$ws.Range("D2:D51").Style = "Normal"
